# Project version numbers.xlsx - update ready for release:
#  - add a "14.0.0" release column (N) with the latest version per project
#  - add a new row documenting the Tardigrade.Framework.MailKit project
#  - freeze the first column and leave the selection on the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "14.0.0" column (N) for every existing project row first.
$ws.Range("N1").Value = "14.0.0"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N2").Value = "12.0.0"
$ws.Range("N3").Value = "6.0.0"
$ws.Range("N4").Value = "5.0.0"
$ws.Range("N5").Value = "3.0.0"
$ws.Range("N6").Value = "8.0.0"
$ws.Range("N7").Value = "12.0.0"
$ws.Range("N8").Value = "11.0.0"
$ws.Range("N9").Value = "2.0.0"
$ws.Range("N10").Value = "4.0.0"
$ws.Range("N11").Value = "2.0.0"

# Insert a new row (9) for Tardigrade.Framework.MailKit, ahead of RestSharp,
# which only has its "14.0.0" version filled in so far.
$ws.Rows.Item(9).Insert() | Out-Null
$ws.Range("A9").Value = "Tardigrade.Framework.MailKit"
$ws.Range("N9").Value = "1.0.0"

# Freeze the first column and restore the active selection.
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A10").Select() | Out-Null
